$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 176

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 163

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 157

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 154

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 132

$ws.Range("B7").Value = 103

$ws.Range("B8").Value = 97
